$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: one before column F ("Fecha Solicitud"),
# one before the (then-shifted) "Disponibilidad" column.
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(9).Insert()

# Header row (row 6) relabeling / new content
$ws.Range("C6").Value = "Cliente"
$ws.Range("D6").Value = "Tomo"
$ws.Range("E6").Value = "Usuario"
$ws.Range("F6").Value = "Movimiento"
$ws.Range("I6").Value = "Fecha Entrega"

# Drop the explicit row height on row 6 so it goes back to the sheet default.
$ws.Rows.Item(6).AutoFit()

# Clear the stray E3 cell that no longer carries formatting after the insert.
$ws.Range("E3").Clear()
